$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "falling" sensor data in columns A:C (rows 2-21) is shifted up by one
# row (the first reading is dropped) and a new final reading is appended
# at the bottom (row 21), repairing/extending the recorded data series.

$values = @(
    @(-0.05105815259351689, -0.02509637922048548, -0.0663297846913338),
    @(-0.1197514058578581, -0.3090105539276482, 0.06041020864532091),
    @(-0.4005748778581615, -0.7805985297475537, 0.1577123148100716),
    @(-0.6522004490806942, -1.329693669364566, 0.1308778794038864),
    @(-0.6102398293358944, -1.220246967815219, 0.4809618578070668),
    @(-0.04692753723689685, -0.2711587122508481, 1.261436768940515),
    @(0.2368920927955998, 2.034886604263666, 0.6732607796078658),
    @(-0.7203119397163371, 3.954537868499759, -2.220546166102086),
    @(-0.9733701603753284, 2.440581185477122, -3.891320841653013),
    @(0.4988514525549754, 0.8867653551555867, -1.71078631139937),
    @(0.1264272814705268, -2.447191684019, -0.4622068021978656),
    @(-0.7295694393771054, -0.004014266388761123, 0.09374600010258789),
    @(-0.1179624412740978, 1.259269575277969, 0.1631955632141657),
    @(-0.1617193005624273, -0.4078179995218952, 0.2281727109636559),
    @(-0.124921940267086, 1.249372124671936, 1.016479730606079),
    @(-0.4928955077415405, -0.4476696934018742, -0.9833766732896989),
    @(-0.2805471434479678, 0.4230750912711692, -0.2188279224293611),
    @(-0.2866266923291341, 0.4182899764605935, -0.003992439912898826),
    @(-0.03713915026968551, 0.07855436143775912, 0.07685266648020034),
    @(0.03008511281084441, 0.0510363349070151, -0.07897615255344478)
)

$row = 2
foreach ($entry in $values) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
